$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $val
}

$changes = @(
    @{Row=2; Col=4; Val='42.991.23'},
    @{Row=2; Col=5; Val='  +1.10%  '},
    @{Row=3; Col=4; Val='2.302.85'},
    @{Row=4; Col=4; Val='1.01'},
    @{Row=4; Col=5; Val='  +0.95%  '},
    @{Row=5; Col=4; Val='309.23'},
    @{Row=5; Col=5; Val='  -2.40%  '},
    @{Row=6; Col=4; Val='105.19'},
    @{Row=6; Col=5; Val='  +1.12%  '},
    @{Row=7; Col=5; Val='  -0.47%  '},
    @{Row=8; Col=5; Val='  -0.18%  '},
    @{Row=9; Col=5; Val='  -0.06%  '},
    @{Row=10; Col=4; Val='39.78'},
    @{Row=10; Col=5; Val='  -0.16%  '},
    @{Row=11; Col=4; Val='0.0908'},
    @{Row=11; Col=5; Val='  +0.34%  '},
    @{Row=12; Col=4; Val='8.28'},
    @{Row=12; Col=5; Val='  -2.81%  '},
    @{Row=13; Col=5; Val='  -0.12%  '},
    @{Row=14; Col=4; Val='0.986'},
    @{Row=14; Col=5; Val='  -0.71%  '},
    @{Row=15; Col=4; Val='15.30'},
    @{Row=15; Col=5; Val='  -0.60%  '},
    @{Row=16; Col=4; Val='2.652.26'},
    @{Row=16; Col=5; Val='  +0.03%  '},
    @{Row=17; Col=4; Val='2.303.76'},
    @{Row=17; Col=5; Val='  -0.14%  '},
    @{Row=18; Col=4; Val='42.651.40'},
    @{Row=18; Col=5; Val='  +0.06%  '},
    @{Row=19; Col=4; Val='7.33'},
    @{Row=19; Col=5; Val='  -3.82%  '},
    @{Row=20; Col=4; Val='13.83'},
    @{Row=20; Col=5; Val='  +0.84%  '},
    @{Row=21; Col=5; Val='  -1.09%  '},
    @{Row=22; Col=4; Val='73.43'},
    @{Row=22; Col=5; Val='  -0.80%  '},
    @{Row=23; Col=5; Val='  -2.87%  '},
    @{Row=24; Col=4; Val='268.07'},
    @{Row=24; Col=5; Val='  +0.51%  '},
    @{Row=25; Col=4; Val='2.23'},
    @{Row=25; Col=5; Val='  -0.38%  '},
    @{Row=26; Col=5; Val='  +0.10%  '},
    @{Row=27; Col=4; Val='7.70'},
    @{Row=27; Col=5; Val='  +16.70%  '},
    @{Row=28; Col=5; Val='  +0.32%  '},
    @{Row=29; Col=4; Val='2.30'},
    @{Row=29; Col=5; Val='  +1.27%  '},
    @{Row=30; Col=4; Val='37.80'},
    @{Row=30; Col=5; Val='  +0.10%  '},
    @{Row=31; Col=4; Val='22.23'},
    @{Row=31; Col=5; Val='  -1.86%  '},
    @{Row=32; Col=4; Val='165.56'},
    @{Row=32; Col=5; Val='  -0.12%  '},
    @{Row=33; Col=5; Val='  -2.09%  '},
    @{Row=34; Col=5; Val='  +6.38%  '},
    @{Row=35; Col=4; Val='0.131'},
    @{Row=35; Col=5; Val='  -0.74%  '},
    @{Row=36; Col=4; Val='0.113'},
    @{Row=36; Col=5; Val='  -0.74%  '},
    @{Row=37; Col=4; Val='4.63'},
    @{Row=37; Col=5; Val='  +1.03%  '},
    @{Row=38; Col=5; Val='  +0.64%  '},
    @{Row=39; Col=4; Val='2.80'},
    @{Row=39; Col=5; Val='  +1.61%  '},
    @{Row=40; Col=5; Val='  -3.17%  '},
    @{Row=41; Col=4; Val='108.13'},
    @{Row=41; Col=5; Val='  +12.52%  '},
    @{Row=42; Col=5; Val='  -3.47%  '},
    @{Row=43; Col=4; Val='71.47'},
    @{Row=43; Col=5; Val='  +1.32%  '},
    @{Row=44; Col=5; Val='  +0.77%  '},
    @{Row=45; Col=5; Val='  +0.25%  '},
    @{Row=46; Col=4; Val='12.25'},
    @{Row=46; Col=5; Val='  -1.57%  '},
    @{Row=47; Col=2; Val='Aave'},
    @{Row=47; Col=3; Val='https://coinranking.com/coin/ixgUfzmLR+aave-aave'},
    @{Row=47; Col=4; Val='111.86'},
    @{Row=47; Col=5; Val='  -5.16%  '},
    @{Row=48; Col=2; Val='Maker'},
    @{Row=48; Col=3; Val='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'},
    @{Row=48; Col=4; Val='1.696.97'},
    @{Row=48; Col=5; Val='  +2.71%  '},
    @{Row=49; Col=4; Val='75.96'},
    @{Row=49; Col=5; Val='  -5.72%  '},
    @{Row=50; Col=5; Val='  -0.42%  '},
    @{Row=51; Col=4; Val='5.17'},
    @{Row=51; Col=5; Val='  -2.30%  '},
)

foreach ($item in $changes) {
    $cell = $ws.Cells.Item($item.Row, $item.Col)
    Set-TextValue $cell $item.Val
}

Write-Host "Applied $($changes.Count) cell updates"
